$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'304.04"
$ws.Range("E2").Formula = "'3.91%"
$ws.Range("E3").Formula = "'14.15%"
$ws.Range("D4").Formula = "'5.059"
$ws.Range("E4").Formula = "'1.82%"
$ws.Range("D6").Formula = "'2.271"
$ws.Range("E6").Formula = "'-0.81%"
$ws.Range("D7").Formula = "'8.106"
$ws.Range("E7").Formula = "'4.25%"
$ws.Range("D8").Formula = "'3.998"
$ws.Range("E8").Formula = "'6.01%"
$ws.Range("D9").Formula = "'0.9274"
$ws.Range("D10").Formula = "'0.09799"
$ws.Range("E10").Formula = "'3.65%"
$ws.Range("D11").Formula = "'0.1819"
$ws.Range("E11").Formula = "'4.69%"
$ws.Range("D12").Formula = "'0.08703"
$ws.Range("E12").Formula = "'4.48%"
$ws.Range("D13").Formula = "'0.03414"
$ws.Range("E13").Formula = "'4.40%"
$ws.Range("D14").Formula = "'0.09925"
$ws.Range("E14").Formula = "'0.19%"
$ws.Range("D15").Formula = "'0.001492"
$ws.Range("E15").Formula = "'-0.29%"
$ws.Range("D16").Formula = "'0.005675"
$ws.Range("E16").Formula = "'-1.30%"
$ws.Range("D17").Formula = "'3.487"
$ws.Range("E17").Formula = "'0.53%"
$ws.Range("E18").Formula = "'-1.91%"
$ws.Range("E19").Formula = "'2.96%"
$ws.Range("D21").Formula = "'4.536"
$ws.Range("E21").Formula = "'10.96%"
$ws.Range("D22").Formula = "'0.2236"
$ws.Range("E22").Formula = "'5.39%"
$ws.Range("D23").Formula = "'0.04674"
$ws.Range("E23").Formula = "'3.29%"
$ws.Range("E24").Formula = "'1.75%"
$ws.Range("E25").Formula = "'4.21%"
$ws.Range("E26").Formula = "'0.15%"
$ws.Range("E27").Formula = "'-20.35%"
$ws.Range("E39").Formula = "'8.34%"
$ws.Range("D40").Formula = "'0.04703"
$ws.Range("E40").Formula = "'2.80%"
$ws.Range("D41").Formula = "'0.007990"
$ws.Range("E41").Formula = "'6.97%"
$ws.Range("E42").Formula = "'4.24%"
$ws.Range("D43").Formula = "'0.008534"
$ws.Range("E43").Formula = "'-13.44%"
$ws.Range("D44").Formula = "'0.002300"
$ws.Range("E44").Formula = "'6.68%"
$ws.Range("D45").Formula = "'0.009137"
$ws.Range("E45").Formula = "'1.19%"
$ws.Range("D46").Formula = "'0.00006139"
$ws.Range("E46").Formula = "'0.73%"
$ws.Range("E47").Formula = "'0.17%"
$ws.Range("D48").Formula = "'5.676"
$ws.Range("E48").Formula = "'113.84%"
$ws.Range("E49").Formula = "'34.72%"
$ws.Range("E50").Formula = "'0.17%"
$ws.Range("D51").Formula = "'0.0002001"
$ws.Range("E51").Formula = "'0.17%"
